# Update Name of Algo
# Apply updated values to the result_data_RandomForest sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "A8"  = -22.33580000000002
    "A10" = -22.12950000000001
    "A12" = -21.509
    "A18" = -22.14690000000001
    "A37" = -19.8205
    "A55" = -22.18789999999998
    "A68" = -21.49719999999999
    "A77" = -20.41679999999999
    "A78" = -19.83779999999998
    "A81" = -21.8493
    "A82" = -21.87469999999999

    "E7"   = 14.9419
    "E15"  = 15.9311
    "E18"  = 17.87850000000001
    "E20"  = 15.88129999999998
    "E29"  = 17.04040000000001
    "E30"  = 15.65999999999999
    "E31"  = 16.1569
    "E40"  = 17.12410000000001
    "E50"  = 16.29959999999999
    "E68"  = 16.9995
    "E76"  = 16.22169999999999
    "E87"  = 16.1419
    "E88"  = 16.3683
    "E96"  = 16.18609999999998
    "E98"  = 15.414
    "E101" = 16.84420000000001
    "E102" = 16.7355
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$wb.Save()
